# Tripadvisor New Orleans shard 54 — update:
#  1. Reorder worksheets so "review_info" comes before "hotel_info".
#  2. Insert a new "State" column into "hotel_info" right after "Hotel_Name",
#     populated with "Louisiana" for the existing hotel row.

$wb = $excel.ActiveWorkbook

# --- 1. Move "review_info" to be the first sheet (swaps tab order with "hotel_info") ---
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wb.Worksheets.Item(1))

# --- 2. Insert "State" column (column C) into "hotel_info", after "Hotel_Name" (column B) ---
$wsHotel = $wb.Worksheets.Item("hotel_info")
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
